# Mise à jour de l'application
# Adds a new daily attendance column (BY) dated 2025-11-12 (serial 45973)
# to the "Présences" sheet, filling in each player's status for that day.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Write the new column's values first -----------------------------
# (Writing values before copying number formats keeps the workbook's
#  dependency graph correctly "dirtied" so COUNTA/COUNTIF formulas that
#  already span the new column recompute below.)

$ws.Range("BY1").Value = 45973

$ws.Range("BY2").Value  = "P"
$ws.Range("BY3").Value  = "P"
$ws.Range("BY4").Value  = "P"
$ws.Range("BY5").Value  = "B"
$ws.Range("BY6").Value  = "B"
$ws.Range("BY7").Value  = "P"
$ws.Range("BY8").Value  = "P"
$ws.Range("BY9").Value  = "P"
$ws.Range("BY10").Value = "P"
$ws.Range("BY11").Value = "P"
# Row 12's player roster ends at column AX (no BX/BY entry for them).
$ws.Range("BY13").Value = "B"
$ws.Range("BY14").Value = "P"
$ws.Range("BY15").Value = "B"
$ws.Range("BY16").Value = "P"
$ws.Range("BY17").Value = "P"
$ws.Range("BY18").Value = "P"
$ws.Range("BY19").Value = "P"
$ws.Range("BY20").Value = "P"
# Row 21 stays blank for this date (still gets the column's formatting).
$ws.Range("BY22").Value = "P"
$ws.Range("BY23").Value = "P"
$ws.Range("BY24").Value = "P"
$ws.Range("BY25").Value = "P"
$ws.Range("BY26").Value = "P"
$ws.Range("BY27").Value = "P"
$ws.Range("BY28").Value = "P"
$ws.Range("BY29").Value = "RH"

# --- 2. Match the new column's formatting to the rest of the table ------
# Column BX carries the formatting (date header style, centered data
# style) that the new BY column should inherit, row by row. Row 12 is
# excluded since it has no BX cell either.

$ws.Range("BX1:BX11").Copy()
$ws.Range("BY1:BY11").PasteSpecial(-4122)

$ws.Range("BX13:BX29").Copy()
$ws.Range("BY13:BY29").PasteSpecial(-4122)

# --- 3. Recalculate so every dependent COUNTA/COUNTIF formula's cached --
#        value reflects the newly added column.
$excel.CalculateFull()

# --- 4. Leave the cursor where the author last left it. ------------------
$ws.Range("CB7").Select() | Out-Null
